# "add participant to groupe"
# - A2 (old "cindy.leschaud@gmail.com") -> new participant email
# - A3 (old "cindy@designpond.ch", hyperlinked) -> new participant email, hyperlink removed, style reset to Normal
# - A4 (old "info@designpond.ch", hyperlinked) -> cleared out, but the hyperlink-style formatting is left in place
# - selection moves from B7 to L7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlinks on A3 and A4 before touching their content.
$ws.Range("A4").Hyperlinks.Delete() | Out-Null
$ws.Range("A3").Hyperlinks.Delete() | Out-Null

# A3 becomes the new "hello@yahoo.fr" entry, with plain (non-hyperlink) styling.
$ws.Range("A3").Value = "hello@yahoo.fr"
$ws.Range("A3").Style = "Normal"

# A2 becomes the newly added participant.
$ws.Range("A2").Value = "droitformation.web@gmail.com"

# A4's old email is removed, leaving just the empty, still-styled cell.
$ws.Range("A4").ClearContents() | Out-Null

# Move the active selection as recorded in the sheet view.
$ws.Range("L7").Select() | Out-Null
